$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("L1").Value = "at_symbol"
$ws.Range("L2").Value = 259.153
$ws.Range("L2").NumberFormat = "@"
$ws.Range("L2").WrapText = $true

$ws.Range("L2").Select()
